$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "English" "Inglés"
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"
Replace-Text "Brief" "Breve"
Replace-Text "An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io" "An email sent to partners in the target country who have sent their documents for review. Se enviará a través de customer.io"
Replace-Text "Target audience" "Público objetivo"
Replace-Text "Thank you for submitting your documents" "Gracias por enviar tus documentos"
Replace-Text "Hi " "Hola "
Replace-Text "Thank you for providing us with your documents for the upcoming " "Gracias por facilitarnos tus documentos para el próximo "
Replace-Text ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation." ". Basándonos en la información que nos has facilitado, haremos los preparativos necesarios, incluidos el alojamiento y el transporte."
Replace-Text "We’re currently reviewing your documents and will reach out to you if we need anything else. " "Estamos revisando tus documentos y nos pondremos en contacto contigo si necesitamos algo más. "
Replace-Text "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-Text " or " " o "
Replace-Text "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-Text ", at " ", en "
Replace-Text "We look forward to seeing you at " "Esperamos verte en "

$c = $d.Comments.Item(1)
$c.Range.Text = "elija uno de los dos"
